# coordenadas.xlsx - agregadas tarjetas de ferrocarriles
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "railroad" card numbers on the board grid ---
$ws.Range("J4").Value = 17
$ws.Range("D10").Value = 10
$ws.Range("P10").Value = 25
$ws.Range("J16").Value = 2

# --- View state: scroll/zoom/selection as last left by the author ---
$win = $excel.ActiveWindow
$win.Zoom = 85

# Move the selection to P10 (the new railroad card cell), matching the
# author's last saved selection/active cell.
$null = $ws.Range("P10").Select()
